$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new BOM row 38: D5 - SMA Diode used for power through microB, with ground fills added.
# New shared strings must be created in this order to match target indices:
#   120 Diode, 121 digikey URL, 122 VS-MBRA140TRPBF, 123 SMA-DIODE, 124 D5
$ws.Range("A38").Value = "N"
$ws.Range("C38").Value = "Diode"
$ws.Range("G38").Value = "http://www.digikey.com/product-detail/en/VS-MBRA140TRPBF/VS-MBRA140TRPBFCT-ND/2687965"
$ws.Range("D38").Value = "VS-MBRA140TRPBF"
$ws.Range("E38").Value = "SMA-DIODE"
$ws.Range("B38").Value = "D5"
$ws.Range("F38").Value = 1

# Update view: scroll/selection moves to just below the new last row.
$ws.Range("B39").Select()
